$d = $word.ActiveDocument
$r = $d.Content
$found = $r.Find.Execute("reasonable demand for a long time.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$origPara = $r.Paragraphs(1)
$r.Collapse(0)
$r.InsertParagraphAfter()
$newPara = $origPara.Next()
$newRange = $newPara.Range
$paraStart = $newRange.Start

$fullText = "These figures have highlighted to me that .NET development is low in demand, which will likely make suitable employment difficult to find. This may be offset by the fact that there appears to be a shortage of .NET developers in Australia. While this will rectify the odds of landing a position when applying, it is important to consider that not all workplaces/positions which require .NET developer skills will be suitable for me. This means that the already small pool of potential jobs is even smaller by the time I eliminate roles or companies that would not be suitable for me and my lifestyle. A small selection of potential job prospects does not mean the career path is unfeasible, it just means that I may benefit from broadening my horizons. This data has not changed my ideal job, but has lead me to realize how greatly I will benefit from maintaining my leadership and interpersonal skills, as well as my programming skills in a wide array of languages. This would greatly widen my list of potential job opportunities."
$newRange.InsertAfter($fullText)
Write-Host "paraStart=$paraStart totalLen=$($fullText.Length)"

# split offsets (relative to paraStart), in the order computed
$splits = @(61, 100, 109, 129, 914, 971)
foreach ($off in $splits) {
    $segStart = $paraStart + $off
    $segEnd = $d.Content.End
    # use end of paragraph's range instead of whole doc; but using a range to end of new paragraph text should be fine too
    $segRange = $d.Range($segStart, $paraStart + $fullText.Length)
    $segRange.Bold = 1
    $segRange.Bold = 0
}

# verify final text
Write-Host "Final paragraph text = [$($newPara.Range.Text)]"
